# Applies the cryptos.xlsx price/volume/coin-order refresh described by the commit
# "Updated cryptos list on Sat Jul 15 16:50:13 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> column -> new value, for every cell that changed between commits.
# Columns D (Price) and E (Volume(1h)) hold text that looks numeric/percentage
# (e.g. "30.346.68", "  -3.26%  "), so each of those cells is forced to Text
# format before assignment to stop Excel from auto-converting it to a number.
$updates = @(
    @{ Row=2; D='30.346.68'; E='  -3.26%  ' }
    @{ Row=3; D='1.937.73'; E='  -3.20%  ' }
    @{ Row=4; D='1.003'; E='  -0.14%  ' }
    @{ Row=5; D='250.96'; E='  -2.11%  ' }
    @{ Row=6; D='0.7105'; E='  -5.69%  ' }
    @{ Row=7; D='1.002'; E='  -0.25%  ' }
    @{ Row=8; D='0.3306'; E='  -3.36%  ' }
    @{ Row=9; D='27.29'; E='  -1.31%  ' }
    @{ Row=10; D='0.07336'; E='  +2.03%  ' }
    @{ Row=11; D='0.8053'; E='  -3.80%  ' }
    @{ Row=12; D='0.08072'; E='  -1.42%  ' }
    @{ Row=13; D='1.936.27'; E='  -3.31%  ' }
    @{ Row=14; D='5.497'; E='  -2.28%  ' }
    @{ Row=15; D='94.47'; E='  -6.31%  ' }
    @{ Row=16; D='15.10'; E='  -3.48%  ' }
    @{ Row=17; D='30.359.00'; E='  -3.17%  ' }
    @{ Row=18; B='BitcoinCash'; C='https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'; D='253.17'; E='  -5.94%  ' }
    @{ Row=19; B='ShibaInu'; C='https://coinranking.com/coin/xz24e0BjL+shibainu-shib'; D='0.000008215'; E='  -1.61%  ' }
    @{ Row=20; D='5.826'; E='  -4.54%  ' }
    @{ Row=21; D='2.192.17'; E='  -3.22%  ' }
    @{ Row=22; D='1.002'; E='  -0.20%  ' }
    @{ Row=23; E='  -0.24%  ' }
    @{ Row=24; D='7.008'; E='  -1.44%  ' }
    @{ Row=25; D='9.731'; E='  -3.68%  ' }
    @{ Row=26; D='163.73'; E='  -0.27%  ' }
    @{ Row=27; D='2.349'; E='  -1.86%  ' }
    @{ Row=28; D='19.32'; E='  -3.41%  ' }
    @{ Row=29; D='0.1298'; E='  -3.68%  ' }
    @{ Row=30; B='Toncoin'; C='https://coinranking.com/coin/67YlI0K1b+toncoin-ton'; D='1.350'; E='  -2.45%  ' }
    @{ Row=31; B='PancakeSwap'; C='https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'; D='1.550'; E='  -3.57%  ' }
    @{ Row=32; D='4.421'; E='  -5.18%  ' }
    @{ Row=33; D='4.160'; E='  -7.00%  ' }
    @{ Row=34; D='0.05188'; E='  -3.77%  ' }
    @{ Row=35; D='1.269'; E='  -2.56%  ' }
    @{ Row=36; D='0.7469'; E='  -5.67%  ' }
    @{ Row=37; D='2.750'; E='  -1.71%  ' }
    @{ Row=38; E='  -2.30%  ' }
    @{ Row=39; D='2.811'; E='  -3.44%  ' }
    @{ Row=40; D='79.01'; E='  -7.91%  ' }
    @{ Row=41; D='6.418'; E='  -6.20%  ' }
    @{ Row=42; E='  -3.23%  ' }
    @{ Row=43; D='2.013'; E='  -5.80%  ' }
    @{ Row=44; D='0.8474'; E='  -1.51%  ' }
    @{ Row=45; D='1.002'; E='  -0.33%  ' }
    @{ Row=46; D='101.68'; E='  -3.82%  ' }
    @{ Row=47; B='EnergySwap'; C='https://coinranking.com/coin/SbWqqTui-+energyswap-ens'; D='9.712'; E='  -4.42%  ' }
    @{ Row=48; B='Aptos'; C='https://coinranking.com/coin/HGYj5JCv5+aptos-apt'; D='7.454'; E='  -4.59%  ' }
    @{ Row=49; D='36.65'; E='  -2.68%  ' }
    @{ Row=50; D='0.4182'; E='  -4.01%  ' }
    @{ Row=51; D='0.06039' }
)

foreach ($update in $updates) {
    $row = $update.Row
    foreach ($col in 'B','C','D','E') {
        if ($update.ContainsKey($col)) {
            $cell = $ws.Range("$col$row")
            if ($col -eq 'D' -or $col -eq 'E') {
                $cell.NumberFormat = '@'
            }
            $cell.Value = $update[$col]
        }
    }
}
